$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data body (A2:C56) before rewriting in the new order/extent.
$ws.Range("A2:C56").ClearContents()

# Copy the date-column style (bold/centered/bordered, style index from A2) so the
# newly-appended rows (57-72) pick up the same formatting as the existing date cells.
$ws.Cells.Item(2,1).Copy()

$ws.Cells.Item(2,1).Value = "2017-10"
$ws.Cells.Item(2,2).Value = 7.95119330558154
$ws.Cells.Item(2,3).Value = 8.227189319887531

$ws.Cells.Item(3,1).Value = "2017-11"
$ws.Cells.Item(3,2).Value = 7.76016838116631
$ws.Cells.Item(3,3).Value = 8.18473287091286

$ws.Cells.Item(4,1).Value = "2017-12"
$ws.Cells.Item(4,2).Value = 7.91234199223232
$ws.Cells.Item(4,3).Value = 8.16203363102281

$ws.Cells.Item(5,1).Value = "2017-02"
$ws.Cells.Item(5,2).Value = 8.24641672924648
$ws.Cells.Item(5,3).Value = 8.24641672924648

$ws.Cells.Item(6,1).Value = "2017-03"
$ws.Cells.Item(6,2).Value = 8.27670828880656
$ws.Cells.Item(6,3).Value = 8.25651391576651

$ws.Cells.Item(7,1).Value = "2017-04"
$ws.Cells.Item(7,2).Value = 8.06298895452511
$ws.Cells.Item(7,3).Value = 8.208132675456159

$ws.Cells.Item(8,1).Value = "2017-05"
$ws.Cells.Item(8,2).Value = 8.0505028270577
$ws.Cells.Item(8,3).Value = 8.176606705776461

$ws.Cells.Item(9,1).Value = "2017-06"
$ws.Cells.Item(9,2).Value = 8.63045380447098
$ws.Cells.Item(9,3).Value = 8.252247888892221

$ws.Cells.Item(10,1).Value = "2017-07"
$ws.Cells.Item(10,2).Value = 8.255315868958631
$ws.Cells.Item(10,3).Value = 8.252686171758841

$ws.Cells.Item(11,1).Value = "2017-08"
$ws.Cells.Item(11,2).Value = 8.34414340032923
$ws.Cells.Item(11,3).Value = 8.264118325330131

$ws.Cells.Item(12,1).Value = "2017-09"
$ws.Cells.Item(12,2).Value = 8.28127215624615
$ws.Cells.Item(12,3).Value = 8.257855543699289

$ws.Cells.Item(13,1).Value = "2018-10"
$ws.Cells.Item(13,2).Value = 7.1932394622197
$ws.Cells.Item(13,3).Value = 7.77540664911443

$ws.Cells.Item(14,1).Value = "2018-11"
$ws.Cells.Item(14,2).Value = 7.2
$ws.Cells.Item(14,3).Value = 7.7

$ws.Cells.Item(15,1).Value = "2018-12"
$ws.Cells.Item(15,2).Value = 7.3
$ws.Cells.Item(15,3).Value = 7.7

$ws.Cells.Item(16,1).Value = "2018-02"
$ws.Cells.Item(16,2).Value = 7.95118911520607
$ws.Cells.Item(16,3).Value = 7.95118911520607

$ws.Cells.Item(17,1).Value = "2018-03"
$ws.Cells.Item(17,2).Value = 8.316649100842049
$ws.Cells.Item(17,3).Value = 8.07300911041807

$ws.Cells.Item(18,1).Value = "2018-04"
$ws.Cells.Item(18,2).Value = 7.98015495302308
$ws.Cells.Item(18,3).Value = 8.049795571069319

$ws.Cells.Item(19,1).Value = "2018-05"
$ws.Cells.Item(19,2).Value = 8.07712382078844
$ws.Cells.Item(19,3).Value = 8.055261221013151

$ws.Cells.Item(20,1).Value = "2018-06"
$ws.Cells.Item(20,2).Value = 7.95693867082225
$ws.Cells.Item(20,3).Value = 8.038874129314671

$ws.Cells.Item(21,1).Value = "2018-07"
$ws.Cells.Item(21,2).Value = 7.56013807161972
$ws.Cells.Item(21,3).Value = 7.97048326392968

$ws.Cells.Item(22,1).Value = "2018-08"
$ws.Cells.Item(22,2).Value = 7.45036579394636
$ws.Cells.Item(22,3).Value = 7.90546858018176

$ws.Cells.Item(23,1).Value = "2018-09"
$ws.Cells.Item(23,2).Value = 7.31707838747053
$ws.Cells.Item(23,3).Value = 7.84009189210274

$ws.Cells.Item(24,1).Value = "2019-10"
$ws.Cells.Item(24,2).Value = 6.6
$ws.Cells.Item(24,3).Value = 7

$ws.Cells.Item(25,1).Value = "2019-11"
$ws.Cells.Item(25,2).Value = 6.8
$ws.Cells.Item(25,3).Value = 6.9

$ws.Cells.Item(26,1).Value = "2019-12"
$ws.Cells.Item(26,2).Value = 6.8
$ws.Cells.Item(26,3).Value = 6.9

$ws.Cells.Item(27,1).Value = "2019-02"
$ws.Cells.Item(27,2).Value = 7.3
$ws.Cells.Item(27,3).Value = 7.3

$ws.Cells.Item(28,1).Value = "2019-03"
$ws.Cells.Item(28,2).Value = 7.6
$ws.Cells.Item(28,3).Value = 7.4

$ws.Cells.Item(29,1).Value = "2019-04"
$ws.Cells.Item(29,2).Value = 7.4
$ws.Cells.Item(29,3).Value = 7.4

$ws.Cells.Item(30,1).Value = "2019-05"
$ws.Cells.Item(30,2).Value = 7
$ws.Cells.Item(30,3).Value = 7.3

$ws.Cells.Item(31,1).Value = "2019-06"
$ws.Cells.Item(31,2).Value = 7.1
$ws.Cells.Item(31,3).Value = 7.3

$ws.Cells.Item(32,1).Value = "2019-07"
$ws.Cells.Item(32,2).Value = 6.3
$ws.Cells.Item(32,3).Value = 7.1

$ws.Cells.Item(33,1).Value = "2019-08"
$ws.Cells.Item(33,2).Value = 6.4
$ws.Cells.Item(33,3).Value = 7

$ws.Cells.Item(34,1).Value = "2019-09"
$ws.Cells.Item(34,2).Value = 6.7
$ws.Cells.Item(34,3).Value = 7

$ws.Cells.Item(35,1).Value = "2020-10"
$ws.Cells.Item(35,2).Value = 7.4
$ws.Cells.Item(35,3).Value = -1.6

$ws.Cells.Item(36,1).Value = "2020-11"
$ws.Cells.Item(36,2).Value = 8
$ws.Cells.Item(36,3).Value = -0.7

$ws.Cells.Item(37,1).Value = "2020-12"
$ws.Cells.Item(37,2).Value = 7.7
$ws.Cells.Item(37,3).Value = 0

$ws.Cells.Item(38,1).Value = "2020-02"
$ws.Cells.Item(38,3).Value = -13

$ws.Cells.Item(39,1).Value = "2020-03"
$ws.Cells.Item(39,2).Value = -9.1
$ws.Cells.Item(39,3).Value = -11.7

$ws.Cells.Item(40,1).Value = "2020-04"
$ws.Cells.Item(40,2).Value = -4.5
$ws.Cells.Item(40,3).Value = -9.9

$ws.Cells.Item(41,1).Value = "2020-05"
$ws.Cells.Item(41,2).Value = 1
$ws.Cells.Item(41,3).Value = -7.7

$ws.Cells.Item(42,1).Value = "2020-06"
$ws.Cells.Item(42,2).Value = 2.3
$ws.Cells.Item(42,3).Value = -6.1

$ws.Cells.Item(43,1).Value = "2020-07"
$ws.Cells.Item(43,2).Value = 3.5
$ws.Cells.Item(43,3).Value = -4.7

$ws.Cells.Item(44,1).Value = "2020-08"
$ws.Cells.Item(44,2).Value = 4
$ws.Cells.Item(44,3).Value = -3.6

$ws.Cells.Item(45,1).Value = "2020-09"
$ws.Cells.Item(45,2).Value = 5.4
$ws.Cells.Item(45,3).Value = -2.6

$ws.Cells.Item(46,1).Value = "2021-10"
$ws.Cells.Item(46,2).Value = 3.8
$ws.Cells.Item(46,3).Value = 15.1

$ws.Cells.Item(47,1).Value = "2021-11"
$ws.Cells.Item(47,2).Value = 3.1
$ws.Cells.Item(47,3).Value = 14

$ws.Cells.Item(48,1).Value = "2021-12"
$ws.Cells.Item(48,2).Value = 3
$ws.Cells.Item(48,3).Value = 13.1

$ws.Cells.Item(49,1).Value = "2021-02"
$ws.Cells.Item(49,3).Value = 31.1

$ws.Cells.Item(50,1).Value = "2021-03"
$ws.Cells.Item(50,2).Value = 25.3
$ws.Cells.Item(50,3).Value = 29.2

$ws.Cells.Item(51,1).Value = "2021-04"
$ws.Cells.Item(51,2).Value = 18.2
$ws.Cells.Item(51,3).Value = 26.4

$ws.Cells.Item(52,1).Value = "2021-05"
$ws.Cells.Item(52,2).Value = 12.5
$ws.Cells.Item(52,3).Value = 23.6

$ws.Cells.Item(53,1).Value = "2021-06"
$ws.Cells.Item(53,2).Value = 10.9
$ws.Cells.Item(53,3).Value = 21.5

$ws.Cells.Item(54,1).Value = "2021-07"
$ws.Cells.Item(54,2).Value = 7.8
$ws.Cells.Item(54,3).Value = 19.6

$ws.Cells.Item(55,1).Value = "2021-08"
$ws.Cells.Item(55,2).Value = 4.8
$ws.Cells.Item(55,3).Value = 17.7

$ws.Cells.Item(56,1).Value = "2021-09"
$ws.Cells.Item(56,2).Value = 5.2
$ws.Cells.Item(56,3).Value = 16.3

$ws.Cells.Item(57,1).Value = "2022-10"
$ws.Cells.Item(57,2).Value = 0.1
$ws.Cells.Item(57,3).Value = 0.1

$ws.Cells.Item(58,1).Value = "2022-11"
$ws.Cells.Item(58,2).Value = -1.9
$ws.Cells.Item(58,3).Value = -0.1

$ws.Cells.Item(59,1).Value = "2022-12"
$ws.Cells.Item(59,2).Value = -0.8
$ws.Cells.Item(59,3).Value = -0.1

$ws.Cells.Item(60,1).Value = "2022-02"
$ws.Cells.Item(60,3).Value = 4.2

$ws.Cells.Item(61,1).Value = "2022-03"
$ws.Cells.Item(61,2).Value = -0.9
$ws.Cells.Item(61,3).Value = 2.5

$ws.Cells.Item(62,1).Value = "2022-04"
$ws.Cells.Item(62,2).Value = -6.1
$ws.Cells.Item(62,3).Value = 0.3

$ws.Cells.Item(63,1).Value = "2022-05"
$ws.Cells.Item(63,2).Value = -5.1
$ws.Cells.Item(63,3).Value = -0.7

$ws.Cells.Item(64,1).Value = "2022-06"
$ws.Cells.Item(64,2).Value = 1.3
$ws.Cells.Item(64,3).Value = -0.4

$ws.Cells.Item(65,1).Value = "2022-07"
$ws.Cells.Item(65,2).Value = 0.6
$ws.Cells.Item(65,3).Value = -0.3

$ws.Cells.Item(66,1).Value = "2022-08"
$ws.Cells.Item(66,2).Value = 1.8
$ws.Cells.Item(66,3).Value = 0

$ws.Cells.Item(67,1).Value = "2022-09"
$ws.Cells.Item(67,2).Value = 1.3
$ws.Cells.Item(67,3).Value = 0.1

$ws.Cells.Item(68,1).Value = "2023-02"
$ws.Cells.Item(68,3).Value = 5.5

$ws.Cells.Item(69,1).Value = "2023-03"
$ws.Cells.Item(69,2).Value = 9.199999999999999
$ws.Cells.Item(69,3).Value = 6.7

$ws.Cells.Item(70,1).Value = "2023-04"
$ws.Cells.Item(70,2).Value = 13.5
$ws.Cells.Item(70,3).Value = 8.4

$ws.Cells.Item(71,1).Value = "2023-05"
$ws.Cells.Item(71,2).Value = 11.7
$ws.Cells.Item(71,3).Value = 9.1

$ws.Cells.Item(72,1).Value = "2023-06"
$ws.Cells.Item(72,2).Value = 6.8
$ws.Cells.Item(72,3).Value = 8.699999999999999

$ws.Range("A57:A72").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select() | Out-Null
